$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.754.27'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '2.100.24'
$ws.Range("E3").Value = '  +0.33%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.06'
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("E6").Value = '  +0.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.32'
$ws.Range("E7").Value = '  +1.96%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  +1.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0840'
$ws.Range("E10").Value = '  -0.21%  '

$ws.Range("E11").Value = '  -1.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.68'
$ws.Range("E12").Value = '  +5.54%  '

$ws.Range("D13").Value = '2.411.36'
$ws.Range("E13").Value = '  +0.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.04'
$ws.Range("E14").Value = '  -1.42%  '

$ws.Range("E15").Value = '  +3.66%  '

$ws.Range("E16").Value = '  +1.10%  '

$ws.Range("D17").Value = '2.127.63'
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("D18").Value = '38.739.44'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.61'
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").Value = '0.0₃0840'
$ws.Range("E21").Value = '  +0.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.01'
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("E24").Value = '  -3.58%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.87'
$ws.Range("E26").Value = '  +0.89%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.53'
$ws.Range("E27").Value = '  +0.79%  '

$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("E29").Value = '  +3.83%  '

$ws.Range("E30").Value = '  +0.83%  '

$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.54'
$ws.Range("E33").Value = '  +1.20%  '

$ws.Range("E34").Value = '  -0.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.04'
$ws.Range("E35").Value = '  +7.53%  '

$ws.Range("E36").Value = '  +1.66%  '

$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.10'
$ws.Range("E40").Value = '  -2.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.68'
$ws.Range("E41").Value = '  +2.64%  '

$ws.Range("E42").Value = '  +3.05%  '

$ws.Range("D43").Value = '1.528.01'
$ws.Range("E43").Value = '  -1.06%  '

$ws.Range("E44").Value = '  +6.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.89'
$ws.Range("E45").Value = '  +1.97%  '

$ws.Range("E46").Value = '  -0.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0913'
$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.14'
$ws.Range("E48").Value = '  -0.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.06'
$ws.Range("E49").Value = '  +2.79%  '

$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").Value = '2.297.49'
$ws.Range("E51").Value = '  +0.32%  '
